{"js": "// \"Optimisation du systeme de meteo\" -- the rain drops now live only in the\n// window paragraph (\"Global\") instead of spanning the whole room, and the\n// \"_GoBack\" marker follows the last edited field (the Bloc section's \"Id\"\n// line) instead of staying on the old \"X\" line.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// 1) Title \"World\" -> \"Global\"\nconst titleParagraph = paragraphs.items[0];\ntitleParagraph.insertText(\"Global\", Word.InsertLocation.replace);\n\n// 2) Drop the spell-check wrapper (w:proofErr) around \"Gravity\" -- clear the\n//    paragraph and retype the same text so no proofErr markup is re-created.\nconst gravityParagraph = paragraphs.items[1];\ngravityParagraph.clear();\ngravityParagraph.insertText(\"Gravity\", Word.InsertLocation.start);\n\nawait context.sync();\n\n// Re-load so we can reach into the two paragraphs that swap the _GoBack\n// bookmark (paragraph content above did not add/remove paragraphs, so the\n// indices of the later items are unchanged).\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst oldBookmarkParagraph = paragraphs.items[16]; // \"X\" under Room, currently holds _GoBack\nconst newBookmarkParagraph = paragraphs.items[24]; // \"Id\" under Bloc, should receive _GoBack\noldBookmarkParagraph.load(\"text\");\nnewBookmarkParagraph.load(\"text\");\nawait context.sync();\n\n// 3a) Rewrite the old bookmark-holder paragraph without the bookmark.\nconst xWithoutBookmark = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:pPr>\n              <w:pStyle w:val=\"Paragraphedeliste\"/>\n              <w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1\"/></w:numPr>\n              <w:spacing w:after=\"0\"/>\n            </w:pPr>\n            <w:r><w:t>X</w:t></w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\noldBookmarkParagraph.insertOoxml(xWithoutBookmark, Word.InsertLocation.replace);\n\n// 3b) Rewrite the new bookmark-holder paragraph with the bookmark appended\n//     after its run, matching where Word left \"_GoBack\" after the edit.\nconst idWithBookmark = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:pPr>\n              <w:pStyle w:val=\"Paragraphedeliste\"/>\n              <w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1\"/></w:numPr>\n              <w:spacing w:after=\"0\"/>\n            </w:pPr>\n            <w:r><w:t>Id</w:t></w:r>\n            <w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>\n            <w:bookmarkEnd w:id=\"0\"/>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\nnewBookmarkParagraph.insertOoxml(idWithBookmark, Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "# \"Optimisation du systeme de meteo\" -- the rain drops now live only in the\n# window paragraph (\"Global\") instead of spanning the whole room, and the\n# \"_GoBack\" marker follows the last edited field (the Bloc section's \"Id\"\n# line) instead of staying on the old \"X\" line.\n\n$d = $word.ActiveDocument\n\n# 1) Title \"World\" -> \"Global\"\n$titleRange = $d.Paragraphs.Item(1).Range\n$titleRange.Text = \"Global\"\n\n# 2) Drop the spell-check wrapper (w:proofErr) around \"Gravity\". Deleting the\n#    run's text (paragraph mark excluded) and re-inserting fresh OOXML for\n#    the paragraph -- with its own pPr -- leaves no proofErr markup behind.\n$gravityRange = $d.Paragraphs.Item(2).Range\n$gravityRange.MoveEnd(1, -1)\n$gravityRange.Delete()\n$gravityXml = @'\n<?xml version=\"1.0\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:pPr>\n              <w:pStyle w:val=\"Paragraphedeliste\"/>\n              <w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1\"/></w:numPr>\n              <w:spacing w:after=\"0\"/>\n            </w:pPr>\n            <w:r><w:t>Gravity</w:t></w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>\n'@\n$gravityRange.InsertXML($gravityXml)\n\n# 3) Move the \"_GoBack\" bookmark off the old \"X\" line (paragraph 17, under\n#    Room) and onto the final \"Id\" line of the Bloc section (paragraph 25).\n$oldBookmarkRange = $d.Paragraphs.Item(17).Range\n$oldBookmarkRange.MoveEnd(1, -1)\n$oldBookmarkRange.Delete()\n$xWithoutBookmarkXml = @'\n<?xml version=\"1.0\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:pPr>\n              <w:pStyle w:val=\"Paragraphedeliste\"/>\n              <w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1\"/></w:numPr>\n              <w:spacing w:after=\"0\"/>\n            </w:pPr>\n            <w:r><w:t>X</w:t></w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>\n'@\n$oldBookmarkRange.InsertXML($xWithoutBookmarkXml)\n\n$newBookmarkRange = $d.Paragraphs.Item(25).Range\n$newBookmarkRange.MoveEnd(1, -1)\n$newBookmarkRange.Delete()\n$idWithBookmarkXml = @'\n<?xml version=\"1.0\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:pPr>\n              <w:pStyle w:val=\"Paragraphedeliste\"/>\n              <w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1\"/></w:numPr>\n              <w:spacing w:after=\"0\"/>\n            </w:pPr>\n            <w:r><w:t>Id</w:t></w:r>\n            <w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>\n            <w:bookmarkEnd w:id=\"0\"/>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>\n'@\n$newBookmarkRange.InsertXML($idWithBookmarkXml)\n"}
